$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update existing rows (swap/rotate match data; id/Div/Div Original Name/Date columns unchanged) ---
# Row 98
$ws.Range("B98").Value = 6800492
$ws.Range("F98").Value = "Barry Town"
$ws.Range("G98").Value = "Penybont"
$ws.Range("H98").Value = 1
$ws.Range("I98").Value = 1
$ws.Range("J98").Value = "D"
$ws.Range("K98").Value = 3.75
$ws.Range("L98").Value = 3.75
$ws.Range("M98").Value = 1.727
$ws.Range("N98").Value = 3
$ws.Range("O98").Value = 3.6
$ws.Range("P98").Value = 2.05
$ws.Range("Q98").Value = 0.25
$ws.Range("R98").Value = 2
$ws.Range("S98").Value = 1.8
$ws.Range("T98").Value = 2.75
$ws.Range("U98").Value = 1.875
$ws.Range("V98").Value = 1.925
$ws.Range("W98").Value = -1
$ws.Range("X98").Value = 2.6
$ws.Range("Y98").Value = -1
$ws.Range("Z98").Value = 0.5
$ws.Range("AA98").Value = -0.5
$ws.Range("AB98").Value = -1
$ws.Range("AC98").Value = 0.925

# Row 99
$ws.Range("B99").Value = 6800493
$ws.Range("F99").Value = "Newtown"
$ws.Range("G99").Value = "Pontypridd Town"
$ws.Range("H99").Value = 3
$ws.Range("I99").Value = 1
$ws.Range("J99").Value = "H"
$ws.Range("K99").Value = 1.45
$ws.Range("L99").Value = 3.75
$ws.Range("M99").Value = 6.5
$ws.Range("N99").Value = 1.444
$ws.Range("O99").Value = 3.8
$ws.Range("P99").Value = 7.5
$ws.Range("Q99").Value = -1.25
$ws.Range("R99").Value = 2.025
$ws.Range("S99").Value = 1.775
$ws.Range("T99").Value = 2.5
$ws.Range("U99").Value = 2
$ws.Range("V99").Value = 1.8
$ws.Range("W99").Value = 0.444
$ws.Range("X99").Value = -1
$ws.Range("Y99").Value = -1
$ws.Range("Z99").Value = 1.025
$ws.Range("AA99").Value = -1
$ws.Range("AB99").Value = 1
$ws.Range("AC99").Value = -1

# Row 100
$ws.Range("B100").Value = 6800495
$ws.Range("F100").Value = "Haverfordwest County"
$ws.Range("G100").Value = "Cardiff MU"
$ws.Range("H100").Value = 1
$ws.Range("I100").Value = 1
$ws.Range("J100").Value = "D"
$ws.Range("K100").Value = 2.5
$ws.Range("L100").Value = 3.4
$ws.Range("M100").Value = 2.4
$ws.Range("N100").Value = 2.45
$ws.Range("O100").Value = 3.2
$ws.Range("P100").Value = 2.7
$ws.Range("Q100").Value = 0
$ws.Range("R100").Value = 1.775
$ws.Range("S100").Value = 2.025
$ws.Range("T100").Value = 2.5
$ws.Range("U100").Value = 1.95
$ws.Range("V100").Value = 1.85
$ws.Range("W100").Value = -1
$ws.Range("X100").Value = 2.2
$ws.Range("Y100").Value = -1
$ws.Range("Z100").Value = 0
$ws.Range("AA100").Value = 0
$ws.Range("AB100").Value = -1
$ws.Range("AC100").Value = 0.8500000000000001

# Row 103
$ws.Range("B103").Value = 6800497
$ws.Range("F103").Value = "Connahs Quay"
$ws.Range("G103").Value = "Barry Town"
$ws.Range("H103").Value = 7
$ws.Range("I103").Value = 0
$ws.Range("J103").Value = "H"
$ws.Range("K103").Value = 1.2
$ws.Range("L103").Value = 6.5
$ws.Range("M103").Value = 9
$ws.Range("N103").Value = 1.166
$ws.Range("O103").Value = 7.5
$ws.Range("P103").Value = 10
$ws.Range("Q103").Value = -2
$ws.Range("R103").Value = 1.85
$ws.Range("S103").Value = 1.95
$ws.Range("T103").Value = 3.25
$ws.Range("U103").Value = 1.825
$ws.Range("V103").Value = 1.975
$ws.Range("W103").Value = 0.1659999999999999
$ws.Range("X103").Value = -1
$ws.Range("Y103").Value = -1
$ws.Range("Z103").Value = 0.8500000000000001
$ws.Range("AA103").Value = -1
$ws.Range("AB103").Value = 0.825
$ws.Range("AC103").Value = -1

# Row 104
$ws.Range("B104").Value = 6800498
$ws.Range("F104").Value = "Newtown"
$ws.Range("G104").Value = "Colwyn Bay"
$ws.Range("H104").Value = 4
$ws.Range("I104").Value = 2
$ws.Range("J104").Value = "H"
$ws.Range("K104").Value = 1.444
$ws.Range("L104").Value = 4.2
$ws.Range("M104").Value = 6
$ws.Range("N104").Value = 1.4
$ws.Range("O104").Value = 4.2
$ws.Range("P104").Value = 6.5
$ws.Range("Q104").Value = -1.25
$ws.Range("R104").Value = 1.925
$ws.Range("S104").Value = 1.875
$ws.Range("T104").Value = 3
$ws.Range("U104").Value = 1.925
$ws.Range("V104").Value = 1.875
$ws.Range("W104").Value = 0.3999999999999999
$ws.Range("X104").Value = -1
$ws.Range("Y104").Value = -1
$ws.Range("Z104").Value = 0.925
$ws.Range("AA104").Value = -1
$ws.Range("AB104").Value = 0.925
$ws.Range("AC104").Value = -1

# Row 108
$ws.Range("B108").Value = 6800503
$ws.Range("F108").Value = "Connahs Quay"
$ws.Range("G108").Value = "Pontypridd Town"
$ws.Range("H108").Value = 3
$ws.Range("I108").Value = 1
$ws.Range("J108").Value = "H"
$ws.Range("K108").Value = 1.2
$ws.Range("L108").Value = 6
$ws.Range("M108").Value = 9
$ws.Range("N108").Value = 1.142
$ws.Range("O108").Value = 7
$ws.Range("P108").Value = 13
$ws.Range("Q108").Value = -2
$ws.Range("R108").Value = 1.775
$ws.Range("S108").Value = 2.025
$ws.Range("T108").Value = 3.25
$ws.Range("U108").Value = 1.95
$ws.Range("V108").Value = 1.85
$ws.Range("W108").Value = 0.1419999999999999
$ws.Range("X108").Value = -1
$ws.Range("Y108").Value = -1
$ws.Range("Z108").Value = 0
$ws.Range("AA108").Value = 0
$ws.Range("AB108").Value = 0.95
$ws.Range("AC108").Value = -1

# Row 110
$ws.Range("B110").Value = 6800049
$ws.Range("F110").Value = "Caernarfon Town"
$ws.Range("G110").Value = "Haverfordwest County"
$ws.Range("H110").Value = 0
$ws.Range("I110").Value = 1
$ws.Range("J110").Value = "A"
$ws.Range("K110").Value = 2.3
$ws.Range("L110").Value = 3.2
$ws.Range("M110").Value = 2.7
$ws.Range("N110").Value = 1.95
$ws.Range("O110").Value = 3.3
$ws.Range("P110").Value = 3.25
$ws.Range("Q110").Value = -0.5
$ws.Range("R110").Value = 2.025
$ws.Range("S110").Value = 1.775
$ws.Range("T110").Value = 2.75
$ws.Range("U110").Value = 1.8
$ws.Range("V110").Value = 2
$ws.Range("W110").Value = -1
$ws.Range("X110").Value = -1
$ws.Range("Y110").Value = 2.25
$ws.Range("Z110").Value = -1
$ws.Range("AA110").Value = 0.7749999999999999
$ws.Range("AB110").Value = -1
$ws.Range("AC110").Value = 1

# Row 146
$ws.Range("B146").Value = 7721608
$ws.Range("F146").Value = "Barry Town"
$ws.Range("G146").Value = "Haverfordwest County"
$ws.Range("H146").Value = 1
$ws.Range("I146").Value = 1
$ws.Range("J146").Value = "D"
$ws.Range("K146").Value = 2.3
$ws.Range("L146").Value = 3.4
$ws.Range("M146").Value = 2.75
$ws.Range("N146").Value = 2.25
$ws.Range("O146").Value = 3.25
$ws.Range("P146").Value = 2.9
$ws.Range("Q146").Value = -0.25
$ws.Range("R146").Value = 2
$ws.Range("S146").Value = 1.8
$ws.Range("T146").Value = 2.5
$ws.Range("U146").Value = 1.925
$ws.Range("V146").Value = 1.875
$ws.Range("W146").Value = -1
$ws.Range("X146").Value = 2.25
$ws.Range("Y146").Value = -1
$ws.Range("Z146").Value = -0.5
$ws.Range("AA146").Value = 0.4
$ws.Range("AB146").Value = -1
$ws.Range("AC146").Value = 0.875

# Row 147
$ws.Range("B147").Value = 7721586
$ws.Range("F147").Value = "Caernarfon Town"
$ws.Range("G147").Value = "TNS"
$ws.Range("H147").Value = 1
$ws.Range("I147").Value = 8
$ws.Range("J147").Value = "A"
$ws.Range("K147").Value = 11
$ws.Range("L147").Value = 8
$ws.Range("M147").Value = 1.142
$ws.Range("N147").Value = 13
$ws.Range("O147").Value = 7.5
$ws.Range("P147").Value = 1.142
$ws.Range("Q147").Value = 2.25
$ws.Range("R147").Value = 1.95
$ws.Range("S147").Value = 1.85
$ws.Range("T147").Value = 3.5
$ws.Range("U147").Value = 1.8
$ws.Range("V147").Value = 2
$ws.Range("W147").Value = -1
$ws.Range("X147").Value = -1
$ws.Range("Y147").Value = 0.1419999999999999
$ws.Range("Z147").Value = -1
$ws.Range("AA147").Value = 0.8500000000000001
$ws.Range("AB147").Value = 0.8
$ws.Range("AC147").Value = -1

# Row 170
$ws.Range("B170").Value = 7721622
$ws.Range("F170").Value = "Pontypridd Town"
$ws.Range("G170").Value = "Haverfordwest County"
$ws.Range("H170").Value = 0
$ws.Range("I170").Value = 1
$ws.Range("J170").Value = "A"
$ws.Range("K170").Value = 2.375
$ws.Range("L170").Value = 3.25
$ws.Range("M170").Value = 2.6
$ws.Range("N170").Value = 3
$ws.Range("O170").Value = 3.2
$ws.Range("P170").Value = 2.25
$ws.Range("Q170").Value = 0.25
$ws.Range("R170").Value = 1.85
$ws.Range("S170").Value = 1.95
$ws.Range("T170").Value = 2.25
$ws.Range("U170").Value = 1.95
$ws.Range("V170").Value = 1.85
$ws.Range("W170").Value = -1
$ws.Range("X170").Value = -1
$ws.Range("Y170").Value = 1.25
$ws.Range("Z170").Value = -1
$ws.Range("AA170").Value = 0.95
$ws.Range("AB170").Value = -1
$ws.Range("AC170").Value = 0.8500000000000001

# Row 171
$ws.Range("B171").Value = 7721621
$ws.Range("F171").Value = "Colwyn Bay"
$ws.Range("G171").Value = "Aberystwyth"
$ws.Range("H171").Value = 1
$ws.Range("I171").Value = 2
$ws.Range("J171").Value = "A"
$ws.Range("K171").Value = 2.1
$ws.Range("L171").Value = 3.2
$ws.Range("M171").Value = 3.1
$ws.Range("N171").Value = 2.3
$ws.Range("O171").Value = 3.25
$ws.Range("P171").Value = 2.75
$ws.Range("Q171").Value = -0.25
$ws.Range("R171").Value = 2.025
$ws.Range("S171").Value = 1.775
$ws.Range("T171").Value = 2.5
$ws.Range("U171").Value = 2
$ws.Range("V171").Value = 1.8
$ws.Range("W171").Value = -1
$ws.Range("X171").Value = -1
$ws.Range("Y171").Value = 1.75
$ws.Range("Z171").Value = -1
$ws.Range("AA171").Value = 0.7749999999999999
$ws.Range("AB171").Value = 1
$ws.Range("AC171").Value = -1

# Row 172
$ws.Range("B172").Value = 7721620
$ws.Range("F172").Value = "Barry Town"
$ws.Range("G172").Value = "Penybont"
$ws.Range("H172").Value = 0
$ws.Range("I172").Value = 0
$ws.Range("J172").Value = "D"
$ws.Range("K172").Value = 3.1
$ws.Range("L172").Value = 3.1
$ws.Range("M172").Value = 2.15
$ws.Range("N172").Value = 4
$ws.Range("O172").Value = 3.2
$ws.Range("P172").Value = 1.909
$ws.Range("Q172").Value = 0.5
$ws.Range("R172").Value = 1.875
$ws.Range("S172").Value = 1.925
$ws.Range("T172").Value = 2.5
$ws.Range("U172").Value = 1.975
$ws.Range("V172").Value = 1.825
$ws.Range("W172").Value = -1
$ws.Range("X172").Value = 2.2
$ws.Range("Y172").Value = -1
$ws.Range("Z172").Value = 0.875
$ws.Range("AA172").Value = -1
$ws.Range("AB172").Value = -1
$ws.Range("AC172").Value = 0.825

# Row 173
$ws.Range("B173").Value = 7721594
$ws.Range("F173").Value = "Connahs Quay"
$ws.Range("G173").Value = "Newtown"
$ws.Range("H173").Value = 0
$ws.Range("I173").Value = 0
$ws.Range("J173").Value = "D"
$ws.Range("K173").Value = 1.4
$ws.Range("L173").Value = 4.5
$ws.Range("M173").Value = 5.75
$ws.Range("N173").Value = 1.55
$ws.Range("O173").Value = 4.5
$ws.Range("P173").Value = 4.5
$ws.Range("Q173").Value = -1
$ws.Range("R173").Value = 1.9
$ws.Range("S173").Value = 1.9
$ws.Range("T173").Value = 3
$ws.Range("U173").Value = 1.8
$ws.Range("V173").Value = 2
$ws.Range("W173").Value = -1
$ws.Range("X173").Value = 3.5
$ws.Range("Y173").Value = -1
$ws.Range("Z173").Value = -1
$ws.Range("AA173").Value = 0.8999999999999999
$ws.Range("AB173").Value = -1
$ws.Range("AC173").Value = 1

# --- Add new rows 174 and 175 (copy formatting from row 173 id/date cells) ---
# Row 174
$ws.Range("A173").Copy()
$ws.Range("A174").PasteSpecial(-4122)
$ws.Range("E173").Copy()
$ws.Range("E174").PasteSpecial(-4122)
$ws.Range("A174").Value = 172
$ws.Range("B174").Value = 7721623
$ws.Range("C174").Value = "Wales Premier League"
$ws.Range("D174").Value = "Wales Premier League"
$ws.Range("E174").Value = 45387.65625
$ws.Range("F174").Value = "Haverfordwest County"
$ws.Range("G174").Value = "Barry Town"
$ws.Range("K174").Value = 1.95
$ws.Range("L174").Value = 3.25
$ws.Range("M174").Value = 3.4
$ws.Range("N174").Value = 1.909
$ws.Range("O174").Value = 3.25
$ws.Range("P174").Value = 3.5
$ws.Range("Q174").Value = -0.5
$ws.Range("R174").Value = 1.975
$ws.Range("S174").Value = 1.825
$ws.Range("T174").Value = 2.5
$ws.Range("U174").Value = 1.9
$ws.Range("V174").Value = 1.9
$ws.Range("W174").Value = 0
$ws.Range("X174").Value = 0
$ws.Range("Y174").Value = 0
$ws.Range("Z174").Value = 0
$ws.Range("AA174").Value = 0

# Row 175
$ws.Range("A173").Copy()
$ws.Range("A175").PasteSpecial(-4122)
$ws.Range("E173").Copy()
$ws.Range("E175").PasteSpecial(-4122)
$ws.Range("A175").Value = 173
$ws.Range("B175").Value = 7721596
$ws.Range("C175").Value = "Wales Premier League"
$ws.Range("D175").Value = "Wales Premier League"
$ws.Range("E175").Value = 45387.65625
$ws.Range("F175").Value = "Bala Town"
$ws.Range("G175").Value = "Newtown"
$ws.Range("K175").Value = 2
$ws.Range("L175").Value = 3.5
$ws.Range("M175").Value = 3
$ws.Range("N175").Value = 1.909
$ws.Range("O175").Value = 3.5
$ws.Range("P175").Value = 3.1
$ws.Range("Q175").Value = -0.5
$ws.Range("R175").Value = 2
$ws.Range("S175").Value = 1.8
$ws.Range("T175").Value = 2.5
$ws.Range("U175").Value = 1.825
$ws.Range("V175").Value = 1.975
$ws.Range("W175").Value = 0
$ws.Range("X175").Value = 0
$ws.Range("Y175").Value = 0
$ws.Range("Z175").Value = 0
$ws.Range("AA175").Value = 0

$excel.CutCopyMode = 0
Write-Host "Edit complete"